# Update the "想去人数" (interested-people count) figures for a few
# exhibition rows on both the "展览" sheet and the combined "全部类型" sheet.
# These two sheets mirror the same event rows, so both need the update.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 578
    $ws.Range("F3").Value = 3632
    $ws.Range("F5").Value = 702
}
